$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -6
$ws.Range("F6").Value = 1
$ws.Range("F8").Value = -11
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = -1
